# Muharrem - LeftNav Update
# Appends three new test-run log rows (44-46) to the ScenarioStatus sheet:
#   Test name | Result | Browser | Date
# The Date column ("12.12.22") looks like a date to Excel's auto-detection,
# so it is entered as a text formula and then pasted back as a value-only
# (Paste Special -> Values) so it lands as a shared-string text cell instead
# of being converted to a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 44; Name = "Login with valid username and password"; Result = "FAILED"; Browser = "chrome"; Date = "12.12.22" },
    @{ Row = 45; Name = "Login with valid username and password"; Result = "PASSED"; Browser = "chrome"; Date = "12.12.22" },
    @{ Row = 46; Name = "Login with valid username and password"; Result = "PASSED"; Browser = "chrome"; Date = "12.12.22" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    $ws.Cells.Item($r.Row, 2).Value = $r.Result
    $ws.Cells.Item($r.Row, 3).Value = $r.Browser

    # Force the date-like string to be stored as text (shared string), not
    # auto-converted to a date serial number.
    $dateCell = $ws.Cells.Item($r.Row, 4)
    $dateCell.Formula = "=""" + $r.Date + """"
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)
}
